# Append the 3rd benchmark run (rows 19-26) to the "fasta-method-1" sheet,
# mirroring the existing two blocks (rows 2-8 and 10-17) with new timing/
# memory data. All text labels reuse values already present on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fasta-method-1")

$ws.Range("A19").Value = "Num procs"
$ws.Range("B19").Value = "Step"
$ws.Range("C19").Value = "Time since prev"
$ws.Range("D19").Value = "Total time"
$ws.Range("E19").Value = "Memory usage"
$ws.Range("F19").Value = "Virtual memory usage"
$ws.Range("H19").Value = "Num procs"
$ws.Range("I19").Value = "Step"
$ws.Range("J19").Value = "Time since prev"
$ws.Range("K19").Value = "Total time"
$ws.Range("L19").Value = "Memory usage"
$ws.Range("M19").Value = "Virtual memory usage"
$ws.Range("O19").Value = "Num procs"
$ws.Range("P19").Value = "Step"
$ws.Range("Q19").Value = "Time since prev"
$ws.Range("R19").Value = "Total time"
$ws.Range("S19").Value = "Memory usage"
$ws.Range("T19").Value = "Virtual memory usage"
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Index ref fasta"
$ws.Range("C20").Value = 0.0026440620422400001
$ws.Range("D20").Value = 0.0026440620422400001
$ws.Range("E20").Value = 12.9921875
$ws.Range("F20").Value = 143.03515625
$ws.Range("H20").Value = 4
$ws.Range("I20").Value = "Index ref fasta"
$ws.Range("J20").Value = 0.0078430175781200005
$ws.Range("K20").Value = 0.0078430175781200005
$ws.Range("L20").Value = 12.99609375
$ws.Range("M20").Value = 143.046875
$ws.Range("O20").Value = 4
$ws.Range("P20").Value = "Index ref fasta"
$ws.Range("Q20").Value = 0.0028359889984099999
$ws.Range("R20").Value = 0.0028359889984099999
$ws.Range("S20").Value = 13
$ws.Range("T20").Value = 143.046875
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Calcs"
$ws.Range("C21").Value = 0.86020994186400002
$ws.Range("D21").Value = 0.86285400390599998
$ws.Range("E21").Value = 31.12109375
$ws.Range("F21").Value = 161.03125
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = "Calcs"
$ws.Range("J21").Value = 0.49518895149199998
$ws.Range("K21").Value = 0.50303196907000003
$ws.Range("L21").Value = 31.12109375
$ws.Range("M21").Value = 161.0390625
$ws.Range("O21").Value = 4
$ws.Range("P21").Value = "Calcs"
$ws.Range("Q21").Value = 0.74214220047000001
$ws.Range("R21").Value = 0.74497818946799998
$ws.Range("S21").Value = 31.125
$ws.Range("T21").Value = 161.0390625
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Split files"
$ws.Range("C22").Value = 0.86113810539199998
$ws.Range("D22").Value = 0.86378216743500003
$ws.Range("E22").Value = 31.12109375
$ws.Range("F22").Value = 161.03125
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = "Split files"
$ws.Range("J22").Value = 0.49603295326199998
$ws.Range("K22").Value = 0.50387597084000002
$ws.Range("L22").Value = 31.12109375
$ws.Range("M22").Value = 161.0390625
$ws.Range("O22").Value = 4
$ws.Range("P22").Value = "Split files"
$ws.Range("Q22").Value = 0.74326419830299995
$ws.Range("R22").Value = 0.74610018730200001
$ws.Range("S22").Value = 31.125
$ws.Range("T22").Value = 161.0390625
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Merge files"
$ws.Range("C23").Value = 0.0013940334320100001
$ws.Range("D23").Value = 407.21117401100003
$ws.Range("E23").Value = 439.80859375
$ws.Range("F23").Value = 1837.078125
$ws.Range("H23").Value = 4
$ws.Range("I23").Value = "Merge files"
$ws.Range("J23").Value = 0.0016298294067399999
$ws.Range("K23").Value = 1990.8568429899999
$ws.Range("L23").Value = 445.12890625
$ws.Range("M23").Value = 1839.82421875
$ws.Range("O23").Value = 4
$ws.Range("P23").Value = "Merge files"
$ws.Range("Q23").Value = 0.0023519992828400001
$ws.Range("R23").Value = 2009.17144704
$ws.Range("S23").Value = 444.0625
$ws.Range("T23").Value = 1838.81640625
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Adding unmapped"
$ws.Range("C24").Value = 19.5485880375
$ws.Range("D24").Value = 426.75976204900002
$ws.Range("E24").Value = 439.82421875
$ws.Range("F24").Value = 1837.078125
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = "Adding unmapped"
$ws.Range("J24").Value = 19.4055550098
$ws.Range("K24").Value = 2010.2623980000001
$ws.Range("L24").Value = 445.14453125
$ws.Range("M24").Value = 1839.82421875
$ws.Range("O24").Value = 4
$ws.Range("P24").Value = "Adding unmapped"
$ws.Range("Q24").Value = 19.793821096399999
$ws.Range("R24").Value = 2028.96526814
$ws.Range("S24").Value = 444.078125
$ws.Range("T24").Value = 1838.81640625
$ws.Range("A25").Value = 4
$ws.Range("B25").Value = "File 1 unmapped done"
$ws.Range("C25").Value = 231.062651157
$ws.Range("D25").Value = 657.82241320599996
$ws.Range("E25").Value = 439.83984375
$ws.Range("F25").Value = 1837.078125
$ws.Range("H25").Value = 4
$ws.Range("I25").Value = "File 1 unmapped done"
$ws.Range("J25").Value = 237.73836612700001
$ws.Range("K25").Value = 2248.0007641299999
$ws.Range("L25").Value = 445.171875
$ws.Range("M25").Value = 1839.82421875
$ws.Range("O25").Value = 4
$ws.Range("P25").Value = "File 1 unmapped done"
$ws.Range("Q25").Value = 241.98763489699999
$ws.Range("R25").Value = 2270.95290303
$ws.Range("S25").Value = 444.10546875
$ws.Range("T25").Value = 1838.81640625
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "End program"
$ws.Range("C26").Value = 0.0030097961425799998
$ws.Range("D26").Value = 657.82542300199998
$ws.Range("E26").Value = 439.83984375
$ws.Range("F26").Value = 1837.078125
$ws.Range("H26").Value = 4
$ws.Range("I26").Value = "End program"
$ws.Range("J26").Value = 0.0030329227447499999
$ws.Range("K26").Value = 2248.0037970499998
$ws.Range("L26").Value = 445.171875
$ws.Range("M26").Value = 1839.82421875
$ws.Range("O26").Value = 4
$ws.Range("P26").Value = "End program"
$ws.Range("Q26").Value = 0.0027699470519999998
$ws.Range("R26").Value = 2270.9556729800001
$ws.Range("S26").Value = 444.10546875
$ws.Range("T26").Value = 1838.81640625

# Move the active selection to where the new data was entered.
$ws.Range("L24").Select() | Out-Null
